$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting existing rows 54..121 down to 55..122
$ws.Rows.Item(54).Insert()

# Populate the new row 54 with the new data point.
# (Columns A,B,C,E,F,G,H,I,J,K,Q,R,T are identical boilerplate shared with neighboring rows.)
$ws.Range("A54").Value = 10
$ws.Range("B54").Value = "Vega Modelo de Temuco"
$ws.Range("C54").Value = "La Araucanía"
$ws.Range("D54").Value = 45264
$ws.Range("E54").Value = 9
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100107
$ws.Range("H54").Value = "Otros"
$ws.Range("I54").Value = 100107011
$ws.Range("J54").Value = "Tuna"
$ws.Range("K54").Value = "Sin especificar"
$ws.Range("L54").Value = "Especial"
$ws.Range("M54").Value = 50
$ws.Range("N54").Value = 40000
$ws.Range("O54").Value = 40000
$ws.Range("P54").Value = 40000
$ws.Range("Q54").Value = "$/caja 18 kilos"
$ws.Range("R54").Value = "Provincia de Los Andes"
$ws.Range("S54").Value = 2222
$ws.Range("T54").Value = 18
